$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Session 4, 5, & 6 Codes
#
# Row 17 (new): S5 / Match
# Row 18: Assert -> S6 / Match
# Row 19: Parallel Execution, now with full S6 detail columns
# Row 20: Config parameters -> Passing headers in request
# Row 21: Scenario outline with example -> Config parameters (with new detail)
# Row 22: Type / String conversion -> Scenario outline with example
# Row 23: Calling other feature/scenarios -> Type / String conversion
# Row 24 (new): Calling other feature/scenarios
# Row 25 (new): Dynamic data passing in request payload
# Row 26: Call (unchanged)
# ---------------------------------------------------------------------------

# Row 17 - S5 / Match
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "S5"
$ws.Range("C17").Value = "Match"
$ws.Range("D17").Value = "* Understading the LHS and RHS of the syntax`n* Schema validation`n* matching array length`n* match text or binary`n* match karate.lowerCase()"
$ws.Range("E17").Value = "Jun 28 2023"
$ws.Range("F17").Value = "Done"
$ws.Range("G17").Value = "Shailendra, Vishal, Ankit G, Ankit J, Pavan, Ashok"
$ws.Range("H17").Value = "Devyani, Disha, Jayant, Himanshu, Nancy"
$ws.Rows("17:17").RowHeight = 72

# Row 18 - S6 / Match (was "Assert")
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "S6"
$ws.Range("C18").Value = "Match"
$ws.Range("D18").Value = "* match contains only`n* match contains any`n* match each"
$ws.Range("E18").Value = "Jul 4 2023"
$ws.Range("F18").Value = "Done"
$ws.Range("G18").Value = "Jayant, Himanshu, Shailendra, Vishal, Ankit J, Pavan, Ashok"
$ws.Range("H18").Value = "Ankit G, Nancy, Devyani, Disha, Nancy"
$ws.Rows("18:18").RowHeight = 43.2

# Row 19 - S6 / Parallel Execution (text already there, now fleshed out)
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "S6"
$ws.Range("C19").Value = "Parallel Execution"
$ws.Range("D19").Value = "* run scenarios in parallel`n* ignore scenarios`n* set env"
$ws.Range("E19").Value = "Jul 4 2023"
$ws.Range("F19").Value = "Done"
$ws.Rows("19:19").RowHeight = 43.2

# Row 20 - Passing headers in request (was "Config parameters")
$ws.Range("A20").Value = 19
$ws.Range("C20").Value = "Passing headers in request"
$ws.Range("D20").Value = "* Passing headers in request individually using key:value pairs`n* Passing headers in request as json`n* configure headers"
$ws.Rows("20:20").RowHeight = 43.2

# Row 21 - Config parameters (was "Scenario outline with example")
$ws.Range("A21").Value = 20
$ws.Range("C21").Value = "Config parameters"
$ws.Range("D21").Value = "* set timeouts`n* set env`n* set global variable"
$ws.Rows("21:21").RowHeight = 43.2

# Row 22 - Scenario outline with example (was "Type / String conversion")
$ws.Range("C22").Value = "Scenario outline with example"

# Row 23 - Type / String conversion (was "Calling other feature/scenarios")
$ws.Range("C23").Value = "Type / String conversion"

# Row 24 (new) - Calling other feature/scenarios
$ws.Range("C24").Value = "Calling other feature/scenarios"

# Row 25 (new) - Dynamic data passing in request payload
$ws.Range("C25").Value = "Dynamic data passing in request payload"

# Row 26 - Call (unchanged)
$ws.Range("C26").Value = "Call"

# ---------------------------------------------------------------------------
# View state: active selection moved along with the edits (matches the new
# location of the row that was being edited, D16 -> D20)
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
try { $win.Panes.Item(2).ScrollRow = 15 } catch { }
try { $win.ScrollRow = 15 } catch { }
$ws.Range("D20").Select()
